$wb = $excel.ActiveWorkbook

# Update RES installed capacities for bus 2 and bus 3 (rows 2 and 3)
$resSheet = $wb.Worksheets.Item("RES installed")
$resSheet.Range("C2").Value = 5
$resSheet.Range("C3").Value = 5

# Recalculate the whole workbook so all dependent formulas (VLOOKUP / SUM / RANDBETWEEN) update
$excel.CalculateFullRebuild()

# Make "RES installed" the active sheet, matching the authored edit
$resSheet.Activate()
$resSheet.Range("F10").Select()
